$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force number-looking price strings to stay text (matches source inlineStr formatting),
# since these values use "." as a thousands separator / carry significant trailing zeros.
$textCellRefs = @('D4', 'D5', 'D6', 'D9', 'D11', 'D12', 'D13', 'D14', 'D18', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D31', 'D32', 'D33', 'D34', 'D35', 'D38', 'D40', 'D42', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($ref in $textCellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.794.02'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '3.740.07'
$ws.Range('E3').Value = '  -2.44%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = '592.46'
$ws.Range('E5').Value = '  -1.36%  '
$ws.Range('D6').Value = '165.37'
$ws.Range('E6').Value = '  -3.54%  '
$ws.Range('D7').Value = '3.737.51'
$ws.Range('E7').Value = '  -2.42%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = '0.518'
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('E10').Value = '  -4.30%  '
$ws.Range('D11').Value = '6.42'
$ws.Range('E11').Value = '  -1.40%  '
$ws.Range('D12').Value = '0.448'
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('D13').Value = '0.0000263'
$ws.Range('E13').Value = '  -6.52%  '
$ws.Range('D14').Value = '35.80'
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('D15').Value = '4.366.20'
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('D16').Value = '3.733.69'
$ws.Range('E16').Value = '  -3.41%  '
$ws.Range('D17').Value = '67.725.34'
$ws.Range('D18').Value = '18.20'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').Value = '  -5.52%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '10.58'
$ws.Range('E21').Value = '  -2.77%  '
$ws.Range('D22').Value = '463.06'
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').Value = '0.701'
$ws.Range('E23').Value = '  -3.93%  '
$ws.Range('D24').Value = '82.79'
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('E25').Value = '  -14.09%  '
$ws.Range('D26').Value = '2.17'
$ws.Range('E26').Value = '  -4.17%  '
$ws.Range('D27').Value = '11.92'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').Value = '10.16'
$ws.Range('E28').Value = '  -2.75%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = '3.886.15'
$ws.Range('E30').Value = '  -2.76%  '
$ws.Range('D31').Value = '2.87'
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').Value = '7.36'
$ws.Range('E32').Value = '  -4.86%  '
$ws.Range('D33').Value = '29.81'
$ws.Range('E33').Value = '  -3.91%  '
$ws.Range('D34').Value = '2.19'
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('D35').Value = '9.02'
$ws.Range('E35').Value = '  -3.66%  '
$ws.Range('D36').Value = '3.689.12'
$ws.Range('E36').Value = '  -3.26%  '
$ws.Range('E37').Value = '  -3.15%  '
$ws.Range('D38').Value = '3.50'
$ws.Range('E38').Value = '  -10.03%  '
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('D40').Value = '0.993'
$ws.Range('E40').Value = '  -2.08%  '
$ws.Range('E41').Value = '  -3.86%  '
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('E44').Value = '  -3.75%  '
$ws.Range('D45').Value = '8.51'
$ws.Range('E45').Value = '  -2.55%  '
$ws.Range('D46').Value = '1.90'
$ws.Range('E46').Value = '  -3.68%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '45.04'
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '393.49'
$ws.Range('E48').Value = '  -5.83%  '
$ws.Range('D49').Value = '144.16'
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('D50').Value = '0.0345'
$ws.Range('E50').Value = '  -3.91%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '25.03'
$ws.Range('E51').Value = '  -2.96%  '
